# Weekly price-sheet update: a new observation (week of 2023-06-13) is
# inserted as the new row 34, pushing the previously-existing rows 34-51
# down to rows 35-52 (dimension grows from A1:R51 to A1:R52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 34..51 down to 35..52, leaving a blank row 34 to fill in.
$ws.Rows.Item(34).Insert()

# Fill the new row 34 with the new weekly record. Columns A, B, C, E, F,
# G, H, I, N, O, Q, R are constant for every record in this sheet.
$ws.Cells.Item(34, 1).Value = 7
$ws.Cells.Item(34, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(34, 3).Value = "Ñuble"
$ws.Cells.Item(34, 4).Value = 45090
$ws.Cells.Item(34, 5).Value = 16
$ws.Cells.Item(34, 6).Value = 100112043
$ws.Cells.Item(34, 7).Value = "Pepino dulce"
$ws.Cells.Item(34, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 60
$ws.Cells.Item(34, 11).Value = 14000
$ws.Cells.Item(34, 12).Value = 14000
$ws.Cells.Item(34, 13).Value = 14000
$ws.Cells.Item(34, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 16).Value = 778
$ws.Cells.Item(34, 17).Value = 18
$ws.Cells.Item(34, 18).Value = "Hortaliza"
